# Refresh the cryptocurrency snapshot: column D ("Price") and column E
# ("Volume(1h)") get new scraped values. Values are written as literal text
# (matching the original inline-string cells) rather than numbers, so
# formatting such as trailing zeros ("35.30") or thousands-dot grouping
# ("64.513.95") survives exactly as scraped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column-D cells whose new text happens to look like a plain number would
# otherwise be auto-converted by Excel into a numeric value (and lose the
# exact text, e.g. "557.57" -> 557.57 or "35.30" -> 35.3). Prefixing the
# value with an apostrophe forces text entry for those; re-applying the
# "Normal" style afterwards clears the resulting quote-prefix style so the
# cell keeps its original (unstyled) appearance.

$ws.Range('D2').Value = '64.513.95'
$ws.Range('E2').Value = '  +4.59%  '
$ws.Range('D3').Value = '3.090.16'
$ws.Range('E3').Value = '  +2.79%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'557.57"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.55%  '
$ws.Range('D6').Value = "'143.94"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.63%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '3.082.49'
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('D9').Value = "'0.499"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').Value = "'7.11"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +16.44%  '
$ws.Range('E11').Value = '  +3.99%  '
$ws.Range('E12').Value = '  +3.83%  '
$ws.Range('D13').Value = "'0.0000227"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.79%  '
$ws.Range('D14').Value = "'35.30"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.43%  '
$ws.Range('D15').Value = '3.596.87'
$ws.Range('E15').Value = '  +2.91%  '
$ws.Range('D16').Value = '64.573.61'
$ws.Range('E16').Value = '  +4.56%  '
$ws.Range('D17').Value = '3.090.63'
$ws.Range('E17').Value = '  +2.82%  '
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('E19').Value = '  +2.66%  '
$ws.Range('D20').Value = "'481.34"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').Value = '  +4.10%  '
$ws.Range('E22').Value = '  +1.01%  '
$ws.Range('D23').Value = "'7.55"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.73%  '
$ws.Range('D24').Value = "'13.33"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.75%  '
$ws.Range('D25').Value = "'80.85"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'2.77"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.97%  '
$ws.Range('E28').Value = '  +5.71%  '
$ws.Range('E29').Value = '  +7.92%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').Value = "'25.99"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.26%  '
$ws.Range('E32').Value = '  +3.04%  '
$ws.Range('E33').Value = '  +5.35%  '
$ws.Range('D34').Value = "'5.67"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('D35').Value = "'6.20"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.14%  '
$ws.Range('D36').Value = "'54.76"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').Value = "'464.83"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.74%  '
$ws.Range('E38').Value = '  +6.38%  '
$ws.Range('D39').Value = "'0.0824"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.66%  '
$ws.Range('E40').Value = '  +17.65%  '
$ws.Range('D41').Value = '3.003.73'
$ws.Range('E41').Value = '  -4.35%  '
$ws.Range('D42').Value = "'8.26"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.42%  '
$ws.Range('E43').Value = '  -1.40%  '
$ws.Range('D44').Value = "'28.17"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.62%  '
$ws.Range('D45').Value = "'0.257"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.82%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('E47').Value = '  +7.89%  '
$ws.Range('E48').Value = '  +3.93%  '
$ws.Range('D49').Value = "'117.70"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.05%  '
$ws.Range('E50').Value = '  +5.86%  '
$ws.Range('E51').Value = '  +2.45%  '
